$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6100189
$ws.Range("I116").Value = 2302.2593
$ws.Range("J116").Value = 17860398
$ws.Range("K116").Value = 2302.2593
$ws.Range("L116").Value = 17860398
$ws.Range("M116").Value = 1139.7407
$ws.Range("N116").Value = -17867282

$ws.Range("H132").Value = 3789318.5
$ws.Range("I132").Value = 4809297
$ws.Range("J132").Value = 827.1429000000001
$ws.Range("K132").Value = 14427891
$ws.Range("L132").Value = 2481.4287
$ws.Range("M132").Value = -14425361
$ws.Range("N132").Value = -7541.4287

$ws.Range("H135").Value = 1179.1
$ws.Range("I135").Value = 477.8718
$ws.Range("J135").Value = 3665.2727
$ws.Range("K135").Value = 4300.8462
$ws.Range("L135").Value = 32987.4543
$ws.Range("M135").Value = -1765.8462
$ws.Range("N135").Value = -38057.4543

$ws.Range("H137").Value = 3801.3333
$ws.Range("I137").Value = 3196
$ws.Range("J137").Value = 6568.5713
$ws.Range("K137").Value = 9588
$ws.Range("L137").Value = 19705.7139
$ws.Range("M137").Value = -7038
$ws.Range("N137").Value = -24805.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 972.13336
$ws.Range("I97").Value = 972.13336
$ws.Range("K97").Value = 972.13336
$ws.Range("M97").Value = -476.13336

$ws.Range("H132").Value = 2359.6487
$ws.Range("I132").Value = 2769.1428
$ws.Range("J132").Value = 1822.1875
$ws.Range("K132").Value = 8307.428400000001
$ws.Range("L132").Value = 5466.5625
$ws.Range("M132").Value = -5777.428400000001
$ws.Range("N132").Value = -10526.5625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2377.4736
$ws.Range("I20").Value = 1800.7
$ws.Range("J20").Value = 3018.3333
$ws.Range("K20").Value = 1800.7
$ws.Range("L20").Value = 3018.3333
$ws.Range("M20").Value = -1553.7
$ws.Range("N20").Value = -3512.3333

$ws.Range("H105").Value = 3319.2856
$ws.Range("I105").Value = 3266.923
$ws.Range("K105").Value = 3266.923
$ws.Range("M105").Value = -1519.923

$ws.Range("H107").Value = 1642.7273
$ws.Range("I107").Value = 1290.1428
$ws.Range("J107").Value = 2259.75
$ws.Range("K107").Value = 1290.1428
$ws.Range("L107").Value = 2259.75
$ws.Range("M107").Value = 629.8571999999999
$ws.Range("N107").Value = -6099.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34573.438
$ws.Range("I31").Value = 46689.78
$ws.Range("K31").Value = 46689.78
$ws.Range("M31").Value = -46394.78

$ws.Range("H34").Value = 34573.438
$ws.Range("I34").Value = 46689.78
$ws.Range("K34").Value = 46689.78
$ws.Range("M34").Value = -46487.78

$ws.Range("H105").Value = 865
$ws.Range("I105").Value = 847.5
$ws.Range("J105").Value = 900
$ws.Range("K105").Value = 847.5
$ws.Range("L105").Value = 900
$ws.Range("M105").Value = 899.5
$ws.Range("N105").Value = -4394

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 229
$ws.Range("I40").Value = 229
$ws.Range("K40").Value = 916
$ws.Range("M40").Value = -847

$ws.Range("H46").Value = 2097.4546
$ws.Range("I46").Value = 1024.4
$ws.Range("J46").Value = 2991.6667
$ws.Range("K46").Value = 3073.2
$ws.Range("L46").Value = 8975.000100000001
$ws.Range("M46").Value = -2982.2
$ws.Range("N46").Value = -9157.000100000001

$ws.Range("H61").Value = 149.5
$ws.Range("I61").Value = 100
$ws.Range("J61").Value = 298
$ws.Range("K61").Value = 300
$ws.Range("L61").Value = 894
$ws.Range("M61").Value = -85
$ws.Range("N61").Value = -1324

$ws.Range("H64").Value = 3174.6667
$ws.Range("I64").Value = 762
$ws.Range("J64").Value = 8000
$ws.Range("K64").Value = 2286
$ws.Range("L64").Value = 24000
$ws.Range("M64").Value = -2016
$ws.Range("N64").Value = -24540

$ws.Range("H67").Value = 3174.6667
$ws.Range("I67").Value = 762
$ws.Range("J67").Value = 8000
$ws.Range("K67").Value = 2286
$ws.Range("L67").Value = 24000
$ws.Range("M67").Value = -1350
$ws.Range("N67").Value = -25872

$ws.Range("H76").Value = 29000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 29000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 87000
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -87766

$ws.Range("H79").Value = 29000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 29000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 87000
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -89652

$ws.Range("H116").Value = 7257.364
$ws.Range("I116").Value = 7249.75
$ws.Range("J116").Value = 7261.7144
$ws.Range("K116").Value = 21749.25
$ws.Range("L116").Value = 21785.1432
$ws.Range("M116").Value = -18307.25
$ws.Range("N116").Value = -28669.1432

$ws.Range("H118").Value = 1823.44
$ws.Range("I118").Value = 480.66666
$ws.Range("K118").Value = 1441.99998
$ws.Range("M118").Value = -198.9999800000001

$ws.Range("H123").Value = 2767.0527
$ws.Range("I123").Value = 863.2857
$ws.Range("J123").Value = 3877.5833
$ws.Range("K123").Value = 2589.8571
$ws.Range("L123").Value = 11632.7499
$ws.Range("M123").Value = -139.8571000000002
$ws.Range("N123").Value = -16532.7499

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H129").Value = 6174181
$ws.Range("I129").Value = 616.6667
$ws.Range("J129").Value = 9260963
$ws.Range("K129").Value = 1850.0001
$ws.Range("L129").Value = 27782889
$ws.Range("M129").Value = 3149.9999
$ws.Range("N129").Value = -27792889

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3954.0977
$ws.Range("I70").Value = 3969.0881
$ws.Range("J70").Value = 3881.2856
$ws.Range("K70").Value = 3969.0881
$ws.Range("L70").Value = 3881.2856
$ws.Range("M70").Value = -3699.0881
$ws.Range("N70").Value = -4421.2856

$ws.Range("H73").Value = 3954.0977
$ws.Range("I73").Value = 3969.0881
$ws.Range("J73").Value = 3881.2856
$ws.Range("K73").Value = 3969.0881
$ws.Range("L73").Value = 3881.2856
$ws.Range("M73").Value = -3033.0881
$ws.Range("N73").Value = -5753.2856

$ws.Range("H113").Value = 2277.875
$ws.Range("I113").Value = 2430
$ws.Range("K113").Value = 2430
$ws.Range("M113").Value = -260

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1252.9474
$ws.Range("I16").Value = 983.6667
$ws.Range("J16").Value = 1714.5714
$ws.Range("K16").Value = 983.6667
$ws.Range("L16").Value = 1714.5714
$ws.Range("M16").Value = -813.6667
$ws.Range("N16").Value = -2054.5714

$ws.Range("H46").Value = 518
$ws.Range("I46").Value = 509.8
$ws.Range("J46").Value = 600
$ws.Range("K46").Value = 509.8
$ws.Range("L46").Value = 600
$ws.Range("M46").Value = -321.8
$ws.Range("N46").Value = -976

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I113").Value = 231.5
$ws.Range("J113").Value = 833.3333
$ws.Range("K113").Value = 694.5
$ws.Range("L113").Value = 2499.9999
$ws.Range("M113").Value = 1475.5
$ws.Range("N113").Value = -6839.9999

$ws.Range("H132").Value = 970.1177
$ws.Range("I132").Value = 807.6667
$ws.Range("J132").Value = 1360
$ws.Range("K132").Value = 2423.0001
$ws.Range("L132").Value = 4080
$ws.Range("M132").Value = 106.9998999999998
$ws.Range("N132").Value = -9140.0092
